$wb = $excel.ActiveWorkbook

# --- Sheet 3: "magnetite data" --- (update selection first, while it is not
# the sheet we want left active, so activating it temporarily doesn't stick)
$ws3 = $wb.Worksheets.Item("magnetite data")
$ws3.Activate()
$ws3.Range("F2").Select()

# --- Sheet 1: "wustite data" ---
$ws1 = $wb.Worksheets.Item("wustite data")
$ws1.Activate()

# Header text updates (O stoichiometry model: rename existing Genge column,
# add a new one converted to Pa)
$ws1.Range("H1").Value = "p from Genge [dynes cm-2]"
$ws1.Range("I1").Value = "p from Genge [Pa]"

# New column I: p from Genge converted with -1 exponent shift (O stoichiometry)
$ws1.Range("I2").Formula = "=10^((11.3-2.0126*10^4/D2)-1)"
$ws1.Range("I3:I8").Formula = "=10^((11.3-2.0126*10^4/D3)-1)"

# Widen column H, size new column I
$ws1.Columns.Item(8).ColumnWidth = 25.5
$ws1.Columns.Item(9).ColumnWidth = 20.67

# Make "wustite data" the active sheet/tab, with a new selection
$ws1.Range("L7").Select()
